$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.650.57"
$ws.Range("E2").Value = "  +0.94%  "
$ws.Range("D3").Value = "1.564.93"
$ws.Range("E3").Value = "  -0.49%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'210.56"
$ws.Range("E5").Value = "  -0.70%  "
$ws.Range("E6").Value = "  -0.52%  "
$ws.Range("D8").Value = "'24.99"
$ws.Range("E8").Value = "  +5.30%  "
$ws.Range("E9").Value = "  -0.40%  "
$ws.Range("E10").Value = "  -0.30%  "
$ws.Range("D11").Value = "'0.0895"
$ws.Range("E11").Value = "  +0.07%  "
$ws.Range("D12").Value = "1.788.05"
$ws.Range("E12").Value = "  -0.50%  "
$ws.Range("D13").Value = "1.565.83"
$ws.Range("E13").Value = "  -0.60%  "
$ws.Range("D14").Value = "28.663.49"
$ws.Range("E14").Value = "  +1.01%  "
$ws.Range("E15").Value = "  -0.31%  "
$ws.Range("D16").Value = "'3.64"
$ws.Range("E16").Value = "  -1.20%  "
$ws.Range("D17").Value = "'61.43"
$ws.Range("E17").Value = "  -0.16%  "
$ws.Range("D18").Value = "'231.78"
$ws.Range("E18").Value = "  +0.74%  "
$ws.Range("E19").Value = "  -0.70%  "
$ws.Range("E20").Value = "  -1.38%  "
$ws.Range("D21").Value = "'0.999"
$ws.Range("E21").Value = "  -0.12%  "
$ws.Range("D22").Value = "'3.91"
$ws.Range("E22").Value = "  -1.02%  "
$ws.Range("D23").Value = "'8.99"
$ws.Range("E23").Value = "  -0.47%  "
$ws.Range("D24").Value = "'2.12"
$ws.Range("E24").Value = "  +3.61%  "
$ws.Range("D25").Value = "'150.74"
$ws.Range("E25").Value = "  -0.27%  "
$ws.Range("D26").Value = "'14.80"
$ws.Range("E26").Value = "  -0.82%  "
$ws.Range("E27").Value = "  -0.27%  "
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("E29").Value = "  -2.18%  "
$ws.Range("E30").Value = "  -4.46%  "
$ws.Range("E31").Value = "  -1.68%  "
$ws.Range("E32").Value = "  -0.98%  "
$ws.Range("D33").Value = "1.390.10"
$ws.Range("E33").Value = "  +0.37%  "
$ws.Range("E34").Value = "  -4.51%  "
$ws.Range("E35").Value = "  -3.01%  "
$ws.Range("E36").Value = "  -1.94%  "
$ws.Range("E37").Value = "  -2.57%  "
$ws.Range("D38").Value = "'2.66"
$ws.Range("E38").Value = "  +0.59%  "
$ws.Range("E39").Value = "  -0.78%  "
$ws.Range("E40").Value = "  +2.20%  "
$ws.Range("E41").Value = "  -0.35%  "
$ws.Range("D42").Value = "'0.999"
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("E43").Value = "  -1.36%  "
$ws.Range("E44").Value = "  -2.37%  "
$ws.Range("D45").Value = "'63.97"
$ws.Range("E45").Value = "  +2.58%  "
$ws.Range("E46").Value = "  -1.90%  "
$ws.Range("D47").Value = "1.700.75"
$ws.Range("E47").Value = "  -0.42%  "
$ws.Range("E48").Value = "  -5.71%  "
$ws.Range("D49").Value = "'85.38"
$ws.Range("E49").Value = "  +0.29%  "
$ws.Range("D50").Value = "'43.17"
$ws.Range("E50").Value = "  +5.66%  "
$ws.Range("E51").Value = "  +0.57%  "

# Clear quote-prefix formatting introduced by auto-detected numeric text
# so the cell style matches the original (General, no explicit style).
$ws.Range("D5").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
